$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 is the "Da 0 a 3 anni" folder. It currently has sections 1-3 filled
# (columns D..L) but is missing "sezione4" (Controllo del movimento), which
# already exists for the "Da 3 a 5 anni" folder in row 3 (columns M:O).
# Replicate that same "Controllo del movimento" section into row 2.

$ws.Range("M2").Value = "Controllo del movimento"
$ws.Range("N2").Value = "controllo.png"
$ws.Range("O2").Value = "Insieme di esercizi finalizzati al miglioramento del controllo della coordinazione motoria."

# The same description text (shared string) is also used by row 3 ("Da 3 a 5
# anni"), so updating the text of the shared string updates it there too.
$ws.Range("O3").Value = "Insieme di esercizi finalizzati al miglioramento del controllo della coordinazione motoria."

# Update the active selection to reflect where the author ended up working.
$ws.Range("P4").Select()

$wb.Save()
